# Fixed Stimulus Absolute Timestamps
$wb = $excel.ActiveWorkbook

# Rename worksheets (new timestamps in their names)
$wb.Worksheets.Item(1).Name = "GNG_TO-1650477841661266"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778438076682"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778438086703"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778438566685"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778439196901"

# Sheet 1 (GNG) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778416282303.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778416432695.csv"
$ws1.Range("B4").Value = "go_stims-16504778416442668.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778416592631.csv"

# Sheet 2 (NB) updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16504778420186946.csv"
$ws2.Range("B3").Value = "OB-1650477842308696.csv"
$ws2.Range("B4").Value = "ZB-match_2-16504778417827005.csv"
$ws2.Range("B5").Value = "TB-165047784284269.csv"
$ws2.Range("B6").Value = "OB-16504778424096625.csv"
$ws2.Range("B7").Value = "ZB-match_7-1650477841880701.csv"
$ws2.Range("B8").Value = "ZB-match_3-1650477841671702.csv"
$ws2.Range("B9").Value = "TB-16504778437916718.csv"
$ws2.Range("B10").Value = "TB-16504778425886886.csv"

# Sheet 3 (RS) updates
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL) updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778438236687.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778438106701.csv"
$ws4.Range("B4").Value = "MM_stims-1650477843839668.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778438236687.csv"
$ws4.Range("B6").Value = "MM_stims-16504778438557017.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778438406684.csv"

# Sheet 5 (vSAT) updates
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650477843871697.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778438596723.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650477843887701.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778439037025.csv"
